$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 22:35"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1640549
$ws.Range("C4").Value = 19647
$ws.Range("D4").Value = 395359
$ws.Range("E4").Value = 1147719
$ws.Range("G4").Value = 1117
$ws.Range("H4").Value = 97471

# Brasil (row 6)
$ws.Range("B6").Value = 319069
$ws.Range("C6").Value = 8148
$ws.Range("E6").Value = 172568
$ws.Range("G6").Value = 459
$ws.Range("H6").Value = 20541

# Francia (row 10)
$ws.Range("B10").Value = 182219
$ws.Range("C10").Value = 393
$ws.Range("E10").Value = 89721

# India (row 14)
$ws.Range("B14").Value = 124794
$ws.Range("C14").Value = 6568
$ws.Range("D14").Value = 51824
$ws.Range("E14").Value = 69244

# Peru (row 15)
$ws.Range("B15").Value = 111698
$ws.Range("C15").Value = 2929
$ws.Range("D15").Value = 44848
$ws.Range("E15").Value = 63606
$ws.Range("G15").Value = 96
$ws.Range("H15").Value = 3244

# Niger's updated figures overtake Republica de Chipre, so Niger moves to
# row 112 (previously Republica de Chipre) and Republica de Chipre drops to
# row 113 (previously Niger), keeping the table sorted by total cases.
$ws.Range("A112").Value = "Niger"
$ws.Range("B112").Value = 937
$ws.Range("C112").Value = 13
$ws.Range("D112").Value = 764
$ws.Range("E112").Value = 113
$ws.Range("H112").Value = 60

$ws.Range("A113").Value = "Republica de Chipre"
$ws.Range("B113").Value = 927
$ws.Range("C113").Value = 4
$ws.Range("D113").Value = 561
$ws.Range("E113").Value = 349
$ws.Range("H113").Value = 17
